$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "DNI" header/content to cell C2 on Hoja1
$ws.Range("C2").Value = "DNI"

# Leave the sheet selection on G9, matching the cursor position
# at the time the workbook was last saved by the author
$ws.Range("G9").Select() | Out-Null
